$wb = $excel.ActiveWorkbook

# Rename quiz sheets to short numeric codes (e.g., "3.1_Arithmetic_Mean" -> "3.1")
$wb.Worksheets.Item(1).Name = "3.1"
$wb.Worksheets.Item(2).Name = "3.2"
$wb.Worksheets.Item(3).Name = "3.3"
$wb.Worksheets.Item(4).Name = "3.4"

# Make the 4th sheet ("3.4") the active tab with a new selected cell
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()
$ws4.Range("C17").Select()
